# Auto-generated edit script: updates the cryptos list price/volume data
# to match the commit "Updated cryptos list ... with GitHub Actions".
#
# Column D holds price text such as "42.691.34" or "2.64" that must stay
# literal text (the sheet stores these as inline strings, not numbers).
# Cells whose new value looks like a plain decimal (e.g. "2.64", "0.570")
# are pre-formatted as Text ("@") before the write so Excel does not silently
# reinterpret them as numbers (which would also risk dropping significant
# trailing zeros, e.g. "0.570" -> 0.57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.691.34"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "2.549.97"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.83"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.89"
$ws.Range("E6").Value = "  +5.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -1.18%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").Value = "  -1.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("E10").Value = "  +0.44%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("E11").Value = "  -1.94%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.43"
$ws.Range("E12").Value = "  -1.72%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -0.66%  "

# Row 14
$ws.Range("D14").Value = "2.942.09"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.02"
$ws.Range("E15").Value = "  +5.76%  "

# Row 16
$ws.Range("D16").Value = "2.567.63"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("E17").Value = "  -1.97%  "

# Row 18
$ws.Range("D18").Value = "42.725.15"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +0.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  -2.20%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0955"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.25"
$ws.Range("E22").Value = "  -0.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.06"
$ws.Range("E23").Value = "  -3.61%  "

# Row 24
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.56"
$ws.Range("E26").Value = "  -2.21%  "

# Row 27
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.72"
$ws.Range("E28").Value = "  -1.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.36"
$ws.Range("E29").Value = "  -1.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.12"
$ws.Range("E30").Value = "  -2.89%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.64"
$ws.Range("E31").Value = "  +0.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("E32").Value = "  -3.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("E33").Value = "  +11.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0805"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.64"
$ws.Range("E35").Value = "  -2.54%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.25"
$ws.Range("E36").Value = "  -3.67%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.06"
$ws.Range("E37").Value = "  -3.56%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.52"
$ws.Range("E38").Value = "  -4.40%  "

# Row 39
$ws.Range("E39").Value = "  -1.76%  "

# Row 40
$ws.Range("E40").Value = "  -0.74%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.25"
$ws.Range("E41").Value = "  +10.77%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.43"
$ws.Range("E42").Value = "  +1.85%  "

# Row 43
$ws.Range("E43").Value = "  +1.12%  "

# Row 44
$ws.Range("E44").Value = "  +0.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -2.03%  "

# Row 46
$ws.Range("D46").Value = "1.976.59"
$ws.Range("E46").Value = "  -1.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.95"
$ws.Range("E47").Value = "  -1.60%  "

# Row 48
$ws.Range("D48").Value = "2.796.58"
$ws.Range("E48").Value = "  +0.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.49"
$ws.Range("E49").Value = "  -4.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.193"
$ws.Range("E50").Value = "  -0.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.49"
$ws.Range("E51").Value = "  -1.88%  "

